# Actualización automática 2025-11-28 15:30:09
#
# Updates sales figures for ALMEIDA CUATIN JHONATHANN CARLOS across the
# three report sheets: "VENTAS POR GRUPO", "VENTA MENSUAL" and
# "CUMPLIMIENTO MENSUAL".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 9 - COMFALASDI COMPAÑIA FAMILIAR LASCANO DIAZ C. LTDA.
$wsGrupo.Range("H9").Value = 711
$wsGrupo.Range("I9").Value = 388.8
$wsGrupo.Range("M9").Value = 2301.89

# Row 21 - MANCHENO PINO HERVIN SANTIAGO
$wsGrupo.Range("M21").Value = 3681.91

# Row 38 - progress counters ("x de 36")
$wsGrupo.Range("H38").Value = "3 de 36"
$wsGrupo.Range("I38").Value = "2 de 36"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F9").Value = 3401.69
$wsMensual.Range("F21").Value = 3779.22
$wsMensual.Range("F38").Value = 7827.839999999999

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Narrow column F slightly (25 -> 24 characters).
$wsCumpl.Columns.Item(6).ColumnWidth = 23.17

# Row 6 - INODOROS
$wsCumpl.Range("D6").Value = 1445.4
$wsCumpl.Range("E6").Value = -784.8000000000001
$wsCumpl.Range("F6").Value = 2.188010899182562

# Row 7 - LAVABOS
$wsCumpl.Range("D7").Value = 675.9
$wsCumpl.Range("E7").Value = -582
$wsCumpl.Range("F7").Value = 7.198083067092651

# Row 12 - PORCELANATO
$wsCumpl.Range("D12").Value = 4301.08
$wsCumpl.Range("E12").Value = 26912.92
$wsCumpl.Range("F12").Value = 0.1377932978791568

# Row 14 - TOTAL
$wsCumpl.Range("D14").Value = 8797.709999999999
$wsCumpl.Range("E14").Value = 31481.85164865473
$wsCumpl.Range("F14").Value = 0.2184162299664408
